$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "contract_coverage_Chickens" correlation row (row 4) with the
# new "yield_Chickens" correlation row.
$ws.Range("B4").Value = "yield_Chickens"
$ws.Range("C4").Value = "High chicken population density may trigger a sale"

# Update the active selection to C4, matching the saved view state.
$ws.Activate()
$ws.Range("C4").Select()
